# Auto-generated edit script: insert two new quarterly columns (D,E)
# with updated financial data, and apply minor restatements to historical
# quarters further down the row (columns F..J in the new layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTE")

# --- Insert two new columns before column D, shifting old D:K to F:M ---
$ws.Range("D1:E1").EntireColumn.Insert()

# --- Copy number formats/styles from the (now-shifted) F:G columns into the
#     freshly inserted D:E columns so the new cells match the existing look ---
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Populate the two new quarter columns (D = most recent, E = prior) ---
$newQuarterData = @(
    @(7, 43465, 43373),
    @(8, 3750000, 3550000),
    @(9, "NA", "NA"),
    @(10, "NA", "NA"),
    @(12, "NA", "NA"),
    @(13, 0, 0),
    @(14, 19000, 10000),
    @(15, 298000, 273000),
    @(17, 3418000, 3121000),
    @(18, 332000, 429000),
    @(20, -12000, 82000),
    @(21, 618000, 784000),
    @(22, 147000, 142000),
    @(23, 173000, 369000),
    @(24, -23000, 34000),
    @(25, 0, 0),
    @(26, 196000, 335000),
    @(27, 191000, 333000),
    @(28, 0, 0),
    @(29, 0, 0),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 12000, -82000),
    @(33, 191000, 333000),
    @(34, 0, 0),
    @(35, 191000, 333000),
    @(38, 43465, 43373),
    @(41, 71000, 84000),
    @(42, 0, 0),
    @(43, 1897000, 1773000),
    @(44, 811000, 767000),
    @(45, 481000, 479000),
    @(46, 3260000, 3103000),
    @(47, 3463000, 3534000),
    @(48, 21650000, 21612000),
    @(49, 3142000, 3157000),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 4773000, 3886000),
    @(53, 0, 0),
    @(54, 36288000, 35292000),
    @(57, 1329000, 1160000),
    @(58, 2108000, 78000),
    @(59, 1001000, 879000),
    @(60, 4438000, 2117000),
    @(61, 12134000, 13620000),
    @(62, 8999000, 8865000),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 26051000, 25085000),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, 6112000, 6093000),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 10237000, 10207000),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, 191000, 333000),
    @(83, 298000, 273000),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 633000, 614000),
    @(91, -57000, -87000),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -1011000, -956000),
    @(96, -161000, -150000),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, 350000, 360000),
    @(101, 0, 0),
    @(102, -28000, 18000),
)

foreach ($entry in $newQuarterData) {
    $r = $entry[0]
    $dVal = $entry[1]
    $eVal = $entry[2]
    if ("$dVal" -ne "x") { $ws.Cells.Item($r, 4).Value = $dVal }
    if ("$eVal" -ne "x") { $ws.Cells.Item($r, 5).Value = $eVal }
}

# --- A handful of older quarters (now columns F..J) were also restated ---
$restatements = @(
    @(17, 8, 2899000),
    @(17, 9, 2811000),
    @(18, 8, 372000),
    @(18, 9, 434000),
    @(20, 8, -26000),
    @(20, 9, 49000),
    @(24, 8, 39000),
    @(26, 8, 175000),
    @(27, 8, 182000),
    @(29, 8, 105000),
    @(32, 8, 26000),
    @(32, 9, -49000),
    @(91, 6, -69000),
    @(91, 7, -61000),
    @(91, 8, -80000),
    @(91, 9, -65000),
    @(91, 10, -46000),
    @(94, 8, -762000),
    @(94, 9, -554000),
    @(102, 8, 1000),
    @(102, 9, 2000),
)

foreach ($entry in $restatements) {
    $r = $entry[0]
    $c = $entry[1]
    $v = $entry[2]
    $ws.Cells.Item($r, $c).Value = $v
}

